# Add the new "valid_segmentation" column (G) to Sheet1.
#
# Header (row 1) gets the new column label, rows 2-11 are tagged
# "segmentation_01" and rows 12-19 are tagged "segmentation_02" —
# mirrors the two TE-pair blocks already present (rows 2-11 use the
# first set, 12-19 the second).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column widths for the new F/G columns (best-fit, matches authored widths) ---
$ws.Columns.Item(6).ColumnWidth = 18.666666666666668
$ws.Columns.Item(7).ColumnWidth = 18

# --- header (inherits the bold header style already applied to row 1) ---
$ws.Range("G1").Value = "valid_segmentation"

# --- segmentation_01 for rows 2-11 ---
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 7).Value = "segmentation_01"
}

# --- segmentation_02 for rows 12-19 ---
for ($r = 12; $r -le 19; $r++) {
    $ws.Cells.Item($r, 7).Value = "segmentation_02"
}

# --- selection / scroll state to match the saved view ---
$ws.Range("G18:G19").Select()
